# Updated cryptos list on Sat May 25 19:46:21 UTC 2024 with GitHub Actions
# Applies per-cell text updates to the crypto price/volume table on the
# active worksheet. Values in column D that look like plain numbers are
# written through a "force text" dance (temporarily apply a text number
# format, write the value, then restore the cell's original Style) so
# that strings like "6.40" or "0.730" stay text instead of being
# auto-coerced to numeric values (which would drop trailing zeros /
# introduce floating-point artifacts and change the cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Ref='D2'; Val='69.150.24'; Numeric=$false},
    @{Ref='E2'; Val='  +0.30%  '; Numeric=$false},
    @{Ref='D3'; Val='3.757.84'; Numeric=$false},
    @{Ref='E3'; Val='  +0.65%  '; Numeric=$false},
    @{Ref='E4'; Val='  -0.02%  '; Numeric=$false},
    @{Ref='D5'; Val='602.45'; Numeric=$true},
    @{Ref='E5'; Val='  +0.22%  '; Numeric=$false},
    @{Ref='D6'; Val='167.12'; Numeric=$true},
    @{Ref='E6'; Val='  -0.40%  '; Numeric=$false},
    @{Ref='D7'; Val='3.756.06'; Numeric=$false},
    @{Ref='E7'; Val='  +0.57%  '; Numeric=$false},
    @{Ref='E8'; Val='  +0.03%  '; Numeric=$false},
    @{Ref='E9'; Val='  +1.31%  '; Numeric=$false},
    @{Ref='E10'; Val='  +2.77%  '; Numeric=$false},
    @{Ref='D11'; Val='6.40'; Numeric=$true},
    @{Ref='E11'; Val='  +1.80%  '; Numeric=$false},
    @{Ref='D12'; Val='0.459'; Numeric=$true},
    @{Ref='E12'; Val='  -0.03%  '; Numeric=$false},
    @{Ref='D13'; Val='38.08'; Numeric=$true},
    @{Ref='E13'; Val='  -0.41%  '; Numeric=$false},
    @{Ref='E14'; Val='  +1.92%  '; Numeric=$false},
    @{Ref='D15'; Val='4.385.48'; Numeric=$false},
    @{Ref='E15'; Val='  +0.60%  '; Numeric=$false},
    @{Ref='D16'; Val='3.751.85'; Numeric=$false},
    @{Ref='E16'; Val='  +0.48%  '; Numeric=$false},
    @{Ref='D17'; Val='69.160.60'; Numeric=$false},
    @{Ref='E17'; Val='  +0.40%  '; Numeric=$false},
    @{Ref='E18'; Val='  +1.55%  '; Numeric=$false},
    @{Ref='B19'; Val='Chainlink'; Numeric=$false},
    @{Ref='C19'; Val='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; Numeric=$false},
    @{Ref='D19'; Val='17.30'; Numeric=$true},
    @{Ref='E19'; Val='  +0.32%  '; Numeric=$false},
    @{Ref='B20'; Val='TRON'; Numeric=$false},
    @{Ref='C20'; Val='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; Numeric=$false},
    @{Ref='D20'; Val='0.114'; Numeric=$true},
    @{Ref='E20'; Val='  -1.52%  '; Numeric=$false},
    @{Ref='D21'; Val='11.09'; Numeric=$true},
    @{Ref='E21'; Val='  +14.33%  '; Numeric=$false},
    @{Ref='D22'; Val='493.57'; Numeric=$true},
    @{Ref='E22'; Val='  -0.79%  '; Numeric=$false},
    @{Ref='D23'; Val='0.730'; Numeric=$true},
    @{Ref='D24'; Val='0.0000152'; Numeric=$true},
    @{Ref='E24'; Val='  +7.13%  '; Numeric=$false},
    @{Ref='D25'; Val='84.96'; Numeric=$true},
    @{Ref='E25'; Val='  +0.15%  '; Numeric=$false},
    @{Ref='E26'; Val='  +0.13%  '; Numeric=$false},
    @{Ref='D27'; Val='12.31'; Numeric=$true},
    @{Ref='E27'; Val='  +0.49%  '; Numeric=$false},
    @{Ref='D28'; Val='10.10'; Numeric=$true},
    @{Ref='E28'; Val='  -0.02%  '; Numeric=$false},
    @{Ref='E29'; Val='  -0.05%  '; Numeric=$false},
    @{Ref='E30'; Val='  +1.35%  '; Numeric=$false},
    @{Ref='D31'; Val='8.21'; Numeric=$true},
    @{Ref='E31'; Val='  +3.59%  '; Numeric=$false},
    @{Ref='E32'; Val='  +2.59%  '; Numeric=$false},
    @{Ref='D33'; Val='31.60'; Numeric=$true},
    @{Ref='E33'; Val='  -0.29%  '; Numeric=$false},
    @{Ref='D34'; Val='3.903.43'; Numeric=$false},
    @{Ref='E34'; Val='  +0.67%  '; Numeric=$false},
    @{Ref='D35'; Val='3.689.45'; Numeric=$false},
    @{Ref='E35'; Val='  +0.74%  '; Numeric=$false},
    @{Ref='E36'; Val='  -0.15%  '; Numeric=$false},
    @{Ref='E37'; Val='  -0.01%  '; Numeric=$false},
    @{Ref='B38'; Val='Filecoin'; Numeric=$false},
    @{Ref='C38'; Val='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Numeric=$false},
    @{Ref='D38'; Val='5.99'; Numeric=$true},
    @{Ref='E38'; Val='  +3.69%  '; Numeric=$false},
    @{Ref='B39'; Val='Mantle'; Numeric=$false},
    @{Ref='C39'; Val='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; Numeric=$false},
    @{Ref='D39'; Val='1.02'; Numeric=$true},
    @{Ref='E39'; Val='  +1.02%  '; Numeric=$false},
    @{Ref='E40'; Val='  +2.84%  '; Numeric=$false},
    @{Ref='D41'; Val='0.326'; Numeric=$true},
    @{Ref='E41'; Val='  +0.92%  '; Numeric=$false},
    @{Ref='E42'; Val='  +5.46%  '; Numeric=$false},
    @{Ref='D43'; Val='430.19'; Numeric=$true},
    @{Ref='E43'; Val='  -1.41%  '; Numeric=$false},
    @{Ref='E44'; Val='  -0.67%  '; Numeric=$false},
    @{Ref='D45'; Val='1.98'; Numeric=$true},
    @{Ref='E45'; Val='  +0.00%  '; Numeric=$false},
    @{Ref='D46'; Val='8.47'; Numeric=$true},
    @{Ref='E46'; Val='  +1.04%  '; Numeric=$false},
    @{Ref='E47'; Val='  +0.00%  '; Numeric=$false},
    @{Ref='D48'; Val='40.41'; Numeric=$true},
    @{Ref='E48'; Val='  -0.19%  '; Numeric=$false},
    @{Ref='D49'; Val='141.21'; Numeric=$true},
    @{Ref='E49'; Val='  -1.52%  '; Numeric=$false},
    @{Ref='D50'; Val='2.795.64'; Numeric=$false},
    @{Ref='E50'; Val='  +1.92%  '; Numeric=$false},
    @{Ref='D51'; Val='0.0353'; Numeric=$true},
    @{Ref='E51'; Val='  +0.53%  '; Numeric=$false}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Ref)
    if ($u.Numeric) {
        $origStyle = $rng.Style
        $rng.NumberFormat = "@"
        $rng.Value = $u.Val
        $rng.Style = $origStyle
    } else {
        $rng.Value = $u.Val
    }
}
